$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.717.13'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.696.01'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '316.71'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = '0.3958'
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").Value = '1.491'
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").Value = '1.004'
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("D11").Value = '52.66'
$ws.Range("E11").Value = '  -8.97%  '
$ws.Range("D12").Value = '0.08932'
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").Value = '7.272'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '23.53'
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").Value = '8.028'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("D16").Value = '0.00001325'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").Value = '1.697.96'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '100.09'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '0.07040'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '19.67'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").Value = '6.998'
$ws.Range("E21").Value = '  +3.89%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").Value = '24.707.13'
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("D25").Value = '3.271'
$ws.Range("E25").Value = '  +8.78%  '
$ws.Range("D26").Value = '2.365'
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").Value = '22.72'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").Value = '162.29'
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("D29").Value = '136.63'
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("D30").Value = '5.165'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").Value = '7.525'
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("D32").Value = '0.08698'
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").Value = '1.055'
$ws.Range("E33").Value = '  -3.15%  '
$ws.Range("D34").Value = '7.081'
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").Value = '11.40'
$ws.Range("E35").Value = '  +3.70%  '
$ws.Range("D36").Value = '0.2744'
$ws.Range("E36").Value = '  +0.94%  '
$ws.Range("D37").Value = '14.49'
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("D38").Value = '1.882'
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").Value = '0.09234'
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").Value = '0.02729'
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").Value = '0.7664'
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D43").Value = '16.16'
$ws.Range("E43").Value = '  +5.03%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '2.595'
$ws.Range("E44").Value = '  +5.45%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.7179'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").Value = '140.38'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = '91.22'
$ws.Range("E50").Value = '  +5.37%  '
$ws.Range("D51").Value = '0.07979'
$ws.Range("E51").Value = '  -0.18%  '
